$wb = $excel.ActiveWorkbook

# --- Registration sheet: update the sample email address in A2 ---
$reg = $wb.Worksheets.Item("Registration")
$reg.Range("A2").Value = "dowr@test.com"
# Minor column-A width nudge (Excel's bestFit autosize reacted to the new text)
$reg.Columns.Item(1).ColumnWidth = 12.6

# --- test_suite sheet: flip Runmode to "Y" for Login / Parameter / VerifyLoginPage ---
$ts = $wb.Worksheets.Item("test_suite")
$ts.Range("B2").Value = "Y"
$ts.Range("B4").Value = "Y"
$ts.Range("B5").Value = "Y"

# Leave the selection on the Runmode column that was edited, matching the
# active sheet/tab at save time.
$ts.Activate()
$ts.Range("B2:B6").Select() | Out-Null
